# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# values for the f674905b... row on the zh-cn / de-de status sheets, and
# mirrors the de-de handoff timestamp into the Overview summary sheet.

$wb = $excel.ActiveWorkbook

$ovw  = $wb.Sheets("Overview")
$zhcn = $wb.Sheets("zh-cn")
$dede = $wb.Sheets("de-de")

# zh-cn: row 3 corresponds to f674905b-9063-4abe-af0b-b694f7aa8537
$zhcn.Range("H3").Value = "2016-08-21 08:53:40"
$zhcn.Range("K3").Value = "2016-08-21 08:53:57"

# de-de: row 3 corresponds to f674905b-9063-4abe-af0b-b694f7aa8537
$dede.Range("H3").Value = "2016-08-21 08:53:44"
$dede.Range("K3").Value = "2016-08-21 08:54:07"

# Overview: "Latest HO Xliff Generate Date" for the f674905b row now
# reflects the freshly generated de-de handoff timestamp.
$ovw.Range("G3").Value = "2016-08-21 08:53:44"
